# Auto update on 2025-12-23 14:45:07
# Applies numeric corrections to several rows in Sheet1 of kp_data.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 89
$ws.Range("O3").Value = 89
$ws.Range("R3").Value = 0.39

# Row 4
$ws.Range("F4").Value = 85
$ws.Range("N4").Value = 339
$ws.Range("Q4").Value = 2.97

# Row 5
$ws.Range("F5").Value = 149
$ws.Range("N5").Value = 149
$ws.Range("Q5").Value = 2.33

# Row 6
$ws.Range("J6").Value = 396
$ws.Range("N6").Value = 396
$ws.Range("P6").Value = 83
$ws.Range("Q6").Value = 4.77

# Row 7
$ws.Range("J7").Value = 231
$ws.Range("N7").Value = 405
$ws.Range("Q7").Value = 4.26

# Row 10
$ws.Range("G10").Value = 86
$ws.Range("O10").Value = 86
$ws.Range("R10").Value = 0.36

# Row 13
$ws.Range("F13").Value = 120
$ws.Range("N13").Value = 120
$ws.Range("P13").Value = 46
$ws.Range("Q13").Value = 2.61

# Row 15
$ws.Range("F15").Value = 39
$ws.Range("N15").Value = 82
$ws.Range("Q15").Value = 3.04

# Row 17
$ws.Range("F17").Value = 25
$ws.Range("N17").Value = 25
$ws.Range("P17").Value = 7
$ws.Range("Q17").Value = 3.57
